# Add a "% of Q Drop's" column (I) to the grade-distribution sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I.
$ws.Range("I1").Value = "% of Q Drop's"

# New data cells: same "0.00%" text already used in F/G/H for each data
# row. Copy/paste (instead of re-typing the literal) so the cell keeps the
# existing shared-string text entry instead of being reinterpreted as a
# percentage number with a new number-format style.
$ws.Range("F3").Copy()
$ws.Range("I3").PasteSpecial()

$ws.Range("F6").Copy()
$ws.Range("I6").PasteSpecial()

$excel.CutCopyMode = $false
